# Applies the "learning note on English" edit described by the XML diff.
$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Insert the new "vocabulary notes" block before the original first
#    paragraph ("Post-school Qualifications"). We replace paragraph 1's
#    whole range (which covers "Post-school Qualifications" + its
#    paragraph mark) with: a "Words:" paragraph, seven new vocabulary
#    paragraphs, an empty paragraph, and finally a new paragraph holding
#    the original "Post-school Qualifications" text.
# ---------------------------------------------------------------------
$p1 = $d.Paragraphs(1)
$r1 = $p1.Range

$vocabXml = @"
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
  <w:pPr><w:jc w:val="left"/></w:pPr>
  <w:r><w:t>Words:</w:t></w:r>
</w:p>
<w:p>
  <w:pPr><w:jc w:val="left"/></w:pPr>
  <w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>hold a diploma</w:t></w:r>
  <w:r><w:t>(formal)</w:t></w:r>
  <w:r><w:t xml:space="preserve"> </w:t></w:r>
  <w:r><w:t xml:space="preserve"> = have a diploma</w:t></w:r>
</w:p>
<w:p>
  <w:pPr><w:jc w:val="left"/></w:pPr>
  <w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>outnumber</w:t></w:r>
  <w:r><w:t xml:space="preserve"> [v] to be greater in number than sb./sth.</w:t></w:r>
</w:p>
<w:p>
  <w:pPr><w:jc w:val="left"/></w:pPr>
  <w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>counterpart</w:t></w:r>
  <w:r><w:t xml:space="preserve"> [C.] a person or thing that has the same positon or function as sb./sth. else in a different place or situation.</w:t></w:r>
  <w:r><w:br/><w:t>e.g. Men with postgraduate diplomas outnumbered their female counterparts.</w:t></w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pBdr><w:bottom w:val="double" w:sz="6" w:space="1" w:color="auto"/></w:pBdr>
    <w:jc w:val="left"/>
  </w:pPr>
  <w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>constitu</w:t></w:r>
  <w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>t</w:t></w:r>
  <w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>e</w:t></w:r>
  <w:r><w:t xml:space="preserve"> [vt.]  to be the parts that together form something.</w:t></w:r>
  <w:r><w:br/><w:t>e.g. Female workers constitute the majority of the labour force.</w:t></w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pBdr><w:bottom w:val="double" w:sz="6" w:space="1" w:color="auto"/></w:pBdr>
    <w:jc w:val="left"/>
  </w:pPr>
  <w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>respectively</w:t></w:r>
  <w:r><w:t xml:space="preserve"> [adv.] in the same order as the people or things already mentioned.</w:t></w:r>
  <w:r><w:br/><w:t>Julie and mark, aged 17 and 19 respectively.</w:t></w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pBdr><w:bottom w:val="double" w:sz="6" w:space="1" w:color="auto"/></w:pBdr>
    <w:jc w:val="left"/>
  </w:pPr>
  <w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">by contrast </w:t></w:r>
  <w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">:  </w:t></w:r>
  <w:r><w:br/></w:r>
  <w:r><w:t xml:space="preserve">e.g. </w:t></w:r>
  <w:r><w:t xml:space="preserve">By contrast, </w:t></w:r>
  <w:r><w:t>more women held undergraduate diplomas.</w:t></w:r>
  <w:r><w:br/></w:r>
  <w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>contrast</w:t></w:r>
  <w:r><w:t xml:space="preserve"> :  a difference between two or more people or things that you can see clearly when they are compared or put close together.</w:t></w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pBdr><w:bottom w:val="double" w:sz="6" w:space="1" w:color="auto"/></w:pBdr>
    <w:jc w:val="left"/>
  </w:pPr>
  <w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">reach </w:t></w:r>
  <w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>level:</w:t></w:r>
  <w:r><w:br/><w:t>reach undergraduate level</w:t></w:r>
</w:p>
<w:p>
  <w:pPr><w:jc w:val="left"/></w:pPr>
</w:p>
<w:p>
  <w:pPr><w:jc w:val="left"/></w:pPr>
  <w:r><w:t>Post-school Qualifications</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@
$r1.InsertXML($vocabXml)

Write-Output "Paragraph count after vocab insert: $($d.Paragraphs.Count)"
